# SqlServer.xlsx: update the SqlIP value used by the "Property" sheet and
# move the active selection, matching the author's manual edit in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 ("SqlIP" row) changes from 127.0.0.1 -> 192.168.0.24.
# Setting .Value with a new string adds it to the shared-string table and
# re-points the cell at it (matches uniqueCount 14 -> 15, new <si> entry,
# and E2's shared-string index changing to the new entry).
$ws.Range("E2").Value = "192.168.0.24"

# The author's last selection before saving moved from G6 to H6.
$ws.Range("H6").Select()
